$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# Re-style C1 / D1 of the merged header (B1:D1): drop the left border so the
# group reads as a single bordered block (top+bottom on C1, top+bottom+right on D1)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# Re-use the formats built above for the matching columns of the second
# merged header (E1:G1), instead of rebuilding them border-by-border, so the
# same two style entries are shared workbook-wide.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Clear the stray empty inline-string cell
$ws2.Range("G5").ClearContents()
